$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 126; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value2 = $cell.Value2 + 1
}
